$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) values must stay literal text (they already were stored as
# inline strings with thousand-dot formatting / trailing zeros in the source data),
# so force Text number format before assignment to stop Excel from coercing them
# into floating point numbers and dropping the original formatting.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.363.94"
$ws.Range("E2").Value = "  +0.98%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.825.70"
$ws.Range("E3").Value = "  +0.01%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.11"
$ws.Range("E5").Value = "  +0.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4468"
$ws.Range("E7").Value = "  -1.69%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3754"
$ws.Range("E8").Value = "  +0.62%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07485"
$ws.Range("E9").Value = "  +2.20%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8859"
$ws.Range("E10").Value = "  +2.97%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.04"
$ws.Range("E11").Value = "  +0.23%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.832.38"
$ws.Range("E12").Value = "  +0.37%  "
$ws.Range("E13").Value = "  +0.98%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.422"
$ws.Range("E14").Value = "  +1.55%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "93.81"
$ws.Range("E15").Value = "  +1.05%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.07113"
$ws.Range("E16").Value = "  +0.28%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.000"
$ws.Range("E17").Value = "  -0.28%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008775"
$ws.Range("E18").Value = "  -0.60%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9999"
$ws.Range("E19").Value = "  -0.07%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.16"
$ws.Range("E20").Value = "  +1.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.356.15"
$ws.Range("E21").Value = "  +0.87%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.414"
$ws.Range("E22").Value = "  +4.31%  "
$ws.Range("E23").Value = "  -0.65%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.059.78"
$ws.Range("E24").Value = "  +0.55%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.967"
$ws.Range("E25").Value = "  -1.80%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.17"
$ws.Range("E26").Value = "  -0.51%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.308"
$ws.Range("E27").Value = "  +3.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.67"
$ws.Range("E28").Value = "  +0.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.380"
$ws.Range("E29").Value = "  +2.30%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "117.83"
$ws.Range("E30").Value = "  +0.34%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08876"
$ws.Range("E31").Value = "  +0.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7873"
$ws.Range("E32").Value = "  +3.78%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.208"
$ws.Range("E33").Value = "  +1.21%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.610"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.914"
$ws.Range("E35").Value = "  -2.41%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9998"
$ws.Range("E36").Value = "  -0.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.110"
$ws.Range("E37").Value = "  +0.55%  "
$ws.Range("E38").Value = "  +1.29%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05300"
$ws.Range("E39").Value = "  +0.22%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.316"
$ws.Range("E40").Value = "  +1.88%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5353"
$ws.Range("E41").Value = "  -0.30%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.860"
$ws.Range("E42").Value = "  -0.92%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1718"
$ws.Range("E43").Value = "  +0.51%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.320"
$ws.Range("E44").Value = "  +17.38%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.671"
$ws.Range("E45").Value = "  +0.55%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5120"
$ws.Range("E46").Value = "  -1.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.58"
$ws.Range("E47").Value = "  -1.12%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.697"
$ws.Range("E48").Value = "  +1.32%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "105.30"
$ws.Range("E49").Value = "  -0.82%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.9998"
$ws.Range("E50").Value = "  -0.01%  "
$ws.Range("E51").Value = "  +0.67%  "
